# Update the worksheet's division answers to match the newly generated set.
# The data lives in table 1; the five populated rows are Word table rows
# 1, 5, 9, 13 and 17 (1-based), each with 5 columns. Addressing cells by
# (row, column) avoids any ambiguity from duplicate text values
# (e.g. "83÷9=9, 2" and "74÷8=9, 2" each occur twice in the original).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
    1  = @("66÷7=9, 3", "64÷6=10, 4", "99÷8=12, 3", "19÷4=4, 3", "58÷3=19, 1")
    5  = @("64÷3=21, 1", "29÷7=4, 1", "45÷8=5, 5", "40÷3=13, 1", "17÷4=4, 1")
    9  = @("96÷9=10, 6", "54÷2=27, 0", "15÷5=3, 0", "54÷9=6, 0", "82÷3=27, 1")
    13 = @("17÷6=2, 5", "54÷3=18, 0", "66÷2=33, 0", "70÷2=35, 0", "25÷9=2, 7")
    17 = @("68÷4=17, 0", "18÷6=3, 0", "94÷3=31, 1", "53÷3=17, 2", "17÷5=3, 2")
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
